# Commit: "error handling and added a new schema, need to remove _id from
# rendering in the ui next"
#
# The "tickets created" sheet's field-schema table (originally on sheet2)
# moves to the "tickets updated" sheet (sheet1, currently empty) with a
# number of field renames/additions/removals, and the "tickets created"
# sheet is rewritten with a small new schema containing just the internal
# "_id" field (marked for later UI removal per the commit message).

$wb = $excel.ActiveWorkbook
$wsUpdated = $wb.Worksheets.Item("tickets updated")
$wsCreated = $wb.Worksheets.Item("tickets created")

# --- New content for "tickets updated" (sheet1) ---------------------------
$updatedRows = @(
    @('_id','internalField','name','__v','checked','readableField','type'),
    @($null,'email','core',0,$true,'Email Id','string'),
    @($null,'externalId','core',0,$true,'Lead Id','string'),
    @($null,'firstName','core',0,$true,'First Name','string'),
    @($null,'lastName','core',0,$false,'Last Name','string'),
    @($null,'source','core',0,$true,'Source','string'),
    @($null,'amount','core',0,$false,'Amount','number'),
    @($null,'followUp','core',0,$false,'Follow Up','date'),
    @($null,'customerEmail','core',0,$false,'Customer Email','string'),
    @($null,'phoneNumberPrefix','core',0,$true,'Country Code','string'),
    @($null,'phoneNumber','core',0,$false,'Mobile Number','string'),
    @($null,'leadStatus','core',0,$true,'Lead Status','string'),
    @($null,'address','core',0,$true,'Address','string'),
    @($null,'companyName','core',0,$true,'Company','string'),
    @($null,'remarks','core',0,$false,'Remarks','string'),
    @($null,'product','core',0,$true,'Product','string'),
    @($null,'geoLocation','core',0,$false,'Geo Location','geo'),
    @($null,'operationalArea','core',0,$true,'Operational Area','string'),
    @($null,'bucket','core',0,$true,'Bucket','string'),
    @($null,'pinCode','core',0,$false,'Pin Code','string'),
    @($null,'createdAt','core',0,$false,'Created At','date'),
    @($null,'updatedAt','core',0,$false,'Updated At','date')
)

# --- New content for "tickets created" (sheet2) ---------------------------
$createdRows = @(
    @('_id','internalField','name','__v','checked','readableField','type'),
    @($null,'_id','core',0,$false,'CRM Id','string')
)

function Set-SheetData {
    param($ws, $rows)

    # Wipe whatever used range currently exists so stale trailing rows
    # (e.g. sheet2's old 28-row table) don't linger past the new data.
    $used = $ws.UsedRange
    if ($used -ne $null) {
        $used.ClearContents()
    }

    for ($r = 0; $r -lt $rows.Count; $r++) {
        $row = $rows[$r]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $val = $row[$c]
            if ($val -ne $null) {
                $ws.Cells.Item($r + 1, $c + 1).Value = $val
            }
        }
    }
}

Set-SheetData $wsUpdated $updatedRows
Set-SheetData $wsCreated $createdRows

Write-Output "done"
